$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header D1: ProjectID -> ProjectKey
$ws.Range("D1").Value = "ProjectKey"

# Update B2: "Test" -> "test"
$ws.Range("B2").Value = "test"

# Update D2: numeric 10000 -> string "TST"
$ws.Range("D2").Value = "TST"

# Update column widths to fit the new header/value text (column D now holds
# "ProjectKey"/"TST", column F holds the wider "test@test.at" hyperlink text)
$ws.Columns.Item(4).ColumnWidth = 8.83
$ws.Columns.Item(6).ColumnWidth = 27.67

# Update selection to G2
$ws.Range("G2").Select()
